# Generate Report for Handback
#
# The handback transform failed for the e4189c83... file in both the
# zh-cn and de-de locales (filename mismatch: 2z3soyp0.kba vs the
# expected handoff-derived name). Reflect that in the report:
#   - Overview sheet: flip that row's Status from "Ready for handoff" to
#     "Handback transform failed".
#   - Per-locale sheets (zh-cn / de-de): fill in the "Error Detail" cell
#     for that row with the mismatch explanation, and widen the column so
#     the message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the e4189c83-522d-4ca1-996f-024af8ffe2aa.md file.
# Its Status columns (E and F) move from "Ready for handoff" to
# "Handback transform failed" -- and so does every other cell in the
# workbook that shared that same "Ready for handoff" string (the Status
# cell for that row on the zh-cn / de-de sheets), since they all point at
# the same shared-string entry.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# zh-cn sheet: row 3 (same file) gets an Error Detail (column P).
$zhcn.Range("P3").Value = "Handback file name: 2z3soyp0.kba is different with handoff file name: e4189c83-522d-4ca1-996f-024af8ffe2aa.7103b928697fedbe7e0ab1e8a595a48cfb6aaee6.zh-cn."

# de-de sheet: row 3 (same file) gets an Error Detail (column P).
$dede.Range("P3").Value = "Handback file name: 2z3soyp0.kba is different with handoff file name: e4189c83-522d-4ca1-996f-024af8ffe2aa.7103b928697fedbe7e0ab1e8a595a48cfb6aaee6.de-de."

# Widen the Error Detail column (P, the 16th column) on both locale
# sheets to 40 characters so the new message fits. 39.17 is the
# ColumnWidth that round-trips to a stored width of 40 (matching the
# other columns already set that way in this workbook).
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
